$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.2419928825622776
$ws1.Range("C2").Value = 0.05973451327433629
$ws1.Range("D2").Value = 0.9642857142857143
$ws1.Range("E2").Value = 0.1125
$ws1.Range("F2").Value = 0.2393617021276596
$ws1.Range("G2").Value = 0.609375
$ws1.Range("H2").Value = 0.7792602996254683
$ws1.Range("I2").Value = 27
$ws1.Range("J2").Value = 425
$ws1.Range("K2").Value = 109
$ws1.Range("L2").Value = 1

# --- Sheet: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.990909090909091
$ws2.Range("C2").Value = 0.2041198501872659
$ws2.Range("D2").Value = 0.3385093167701863

$ws2.Range("B3").Value = 0.05973451327433629
$ws2.Range("C3").Value = 0.9642857142857143
$ws2.Range("D3").Value = 0.1125

$ws2.Range("B4").Value = 0.2419928825622776
$ws2.Range("C4").Value = 0.2419928825622776
$ws2.Range("D4").Value = 0.2419928825622776
$ws2.Range("E4").Value = 0.2419928825622776

$ws2.Range("B5").Value = 0.5253218020917136
$ws2.Range("C5").Value = 0.5842027822364901
$ws2.Range("D5").Value = 0.2255046583850931

$ws2.Range("B6").Value = 0.9445160514539787
$ws2.Range("C6").Value = 0.2419928825622776
$ws2.Range("D6").Value = 0.327249066112597

# --- Sheet: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 109
$ws3.Range("C2").Value = 425
$ws3.Range("B3").Value = 1
$ws3.Range("C3").Value = 27
